$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("ExcelModuleDemoToDoItem").Name = "DemoToDoItem"
$wb.Worksheets.Item("ExcelModuleDemoProfessionalImp").Name = "ExcelUploadRowHandler4ToDoItem"
